$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '57.123.53'
$ws.Range('E2').Value = '  +4.43%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.363.91'
$ws.Range('E3').Value = '  +3.37%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '519.73'
$ws.Range('E5').Value = '  +3.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '135.21'
$ws.Range('E6').Value = '  +3.78%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.39%  '
$ws.Range('E8').Value = '  +1.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.360.91'
$ws.Range('E9').Value = '  +2.65%  '
$ws.Range('E10').Value = '  +8.08%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.155'
$ws.Range('E11').Value = '  +0.92%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '5.23'
$ws.Range('E12').Value = '  +6.44%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.343'
$ws.Range('E13').Value = '  +1.49%  '
$ws.Range('B14').Value = 'Avalanche'
$ws.Range('C14').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '23.90'
$ws.Range('E14').Value = '  +2.92%  '
$ws.Range('B15').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C15').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.761.20'
$ws.Range('E15').Value = '  +2.51%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '57.007.43'
$ws.Range('E16').Value = '  +4.22%  '
$ws.Range('E17').Value = '  +3.40%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.382.89'
$ws.Range('E18').Value = '  +3.22%  '
$ws.Range('E19').Value = '  +2.19%  '
$ws.Range('E20').Value = '  +2.96%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '322.69'
$ws.Range('E21').Value = '  +5.30%  '
$ws.Range('E22').Value = '  +5.98%  '
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '61.51'
$ws.Range('E24').Value = '  +0.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.993'
$ws.Range('E25').Value = '  -0.23%  '
$ws.Range('E26').Value = '  +6.22%  '
$ws.Range('E27').Value = '  +5.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '170.93'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0₃0741'
$ws.Range('E29').Value = '  +4.66%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.20'
$ws.Range('E30').Value = '  +9.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.30'
$ws.Range('E31').Value = '  +4.16%  '
$ws.Range('E32').Value = '  +3.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.39'
$ws.Range('E33').Value = '  +2.19%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  +0.04%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('E35').Value = '  +0.37%  '
$ws.Range('B36').Value = 'SuiNetwork'
$ws.Range('C36').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.955'
$ws.Range('E36').Value = '  +1.83%  '
$ws.Range('E37').Value = '  +4.92%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.03'
$ws.Range('E38').Value = '  +7.22%  '
$ws.Range('E39').Value = '  +7.10%  '
$ws.Range('E40').Value = '  +3.48%  '
$ws.Range('E41').Value = '  +1.73%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '140.47'
$ws.Range('E42').Value = '  +11.67%  '
$ws.Range('E43').Value = '  +5.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '278.35'
$ws.Range('E44').Value = '  +12.77%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.15'
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0511'
$ws.Range('E46').Value = '  +3.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0932'
$ws.Range('E47').Value = '  +3.55%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.563'
$ws.Range('E48').Value = '  +2.45%  '
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('E50').Value = '  +4.46%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '16.98'
$ws.Range('E51').Value = '  +2.94%  '
